# Insert a new data row at row 356 (pushes existing rows 356-452 down to
# 357-453) and populate the newly inserted row with its values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Rows("356").Insert()

$ws.Range("A356").Value = 6
$ws.Range("B356").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C356").Value = "Metropolitana"
$ws.Range("D356").Value = 44736
$ws.Range("E356").Value = 13
$ws.Range("F356").Value = 100112043
$ws.Range("G356").Value = "Pepino ensalada"
$ws.Range("H356").Value = "Sin especificar"
$ws.Range("I356").Value = "Primera"
$ws.Range("J356").Value = 220
$ws.Range("K356").Value = 15000
$ws.Range("L356").Value = 16000
$ws.Range("M356").Value = 15455
$ws.Range("N356").Value = "`$/caja 60 unidades"
$ws.Range("O356").Value = "Región de Arica y Parinacota"
$ws.Range("P356").Value = 258
$ws.Range("Q356").Value = 60
$ws.Range("R356").Value = "Hortaliza"
